$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "0.999", "70.10") are stored as text and keep trailing/insignificant
# digits exactly as scraped, instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.119.19"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.962.57"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "380.38"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "102.48"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").Value = "36.64"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "3.423.08"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "7.78"
$ws.Range("E14").Value = "  +6.38%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "18.35"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "11.99"
$ws.Range("E16").Value = "  +67.78%  "
$ws.Range("D17").Value = "2.967.44"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").Value = "51.160.79"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "3.11"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "70.10"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("E24").Value = "  +14.26%  "
$ws.Range("D25").Value = "268.05"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").Value = "7.94"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "7.21"
$ws.Range("E27").Value = "  -7.03%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "25.93"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "0.110"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("E32").Value = "  +6.06%  "
$ws.Range("D33").Value = "34.46"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").Value = "51.04"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "0.0437"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +9.78%  "
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("D41").Value = "16.63"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("D42").Value = "2.51"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").Value = "124.67"
$ws.Range("E43").Value = "  +3.43%  "
$ws.Range("D44").Value = "21.62"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").Value = "3.54"
$ws.Range("E45").Value = "  +10.08%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "2.02"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").Value = "0.270"
$ws.Range("E48").Value = "  -4.82%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.050.19"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("E50").Value = "  -5.89%  "
$ws.Range("E51").Value = "  +7.69%  "

# Reset column D style back to the default/Normal style so no stray
# cell-format (style index) diff is introduced; the cells remain text
# because their stored value is already a string.
$ws.Range("D2:D51").Style = "Normal"

